$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append two new log rows (144 and 145) to the feed logs sheet
$ws.Range("A144").Value = 143
$ws.Range("B144").Value = 1
$ws.Range("C144").Value = "2024-06-17 23:13:18"
$ws.Range("D144").Value = 200
$ws.Range("E144").Value = 13

$ws.Range("A145").Value = 144
$ws.Range("B145").Value = 2
$ws.Range("C145").Value = "2024-06-17 23:13:18"
$ws.Range("D145").Value = 200
$ws.Range("E145").Value = 0
